# Updates the cryptocurrency price / 1h-volume table with refreshed market
# data (GitHub Actions scheduled refresh).
#
# Column D ("Price") holds values such as "2.80", "46.704.48" or "0.0₃0985"
# that are stored as plain text in the workbook. If we just assign a string
# like "2.80" to a cell, Excel's automatic type inference will silently turn
# it into the *number* 2.8 (dropping the trailing zero) or mangle values
# such as "301.35" into a floating point number, which does not match the
# source data. To avoid that we briefly force the cell to Text format,
# write the literal string, and then restore the cell's original (default)
# formatting so the saved file doesn't end up with stray per-cell number
# formats that weren't there before.
#
# Column E ("Volume(1h)") holds percentage-like strings such as "  -0.43%  "
# (padded with spaces) which Excel never reinterprets as numbers, so those
# can be written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $cell = $ws.Range($rangeAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

function Set-PlainValue($rangeAddress, $value) {
    $ws.Range($rangeAddress).Value = $value
}

# row, price(D), volume(E)
$rows = @(
    @{ R = 2;  D = "46.704.48";  E = "  -0.43%  " },
    @{ R = 3;  D = "2.276.10";   E = "  -2.28%  " },
    @{ R = 4;  D = $null;        E = "  -0.05%  " },
    @{ R = 5;  D = "301.35";     E = "  -1.58%  " },
    @{ R = 6;  D = "100.03";     E = "  +2.17%  " },
    @{ R = 7;  D = "0.572";      E = "  -1.07%  " },
    @{ R = 8;  D = $null;        E = "  +0.09%  " },
    @{ R = 9;  D = "0.509";      E = "  -5.11%  " },
    @{ R = 10; D = "35.29";      E = "  -1.62%  " },
    @{ R = 11; D = "0.0806";     E = "  -0.13%  " },
    @{ R = 12; D = "7.08";       E = "  -4.96%  " },
    @{ R = 13; D = $null;        E = "  -1.61%  " },
    @{ R = 14; D = "2.623.71";   E = "  -2.23%  " },
    @{ R = 15; D = "2.276.11";   E = "  -2.69%  " },
    @{ R = 16; D = $null;        E = "  -2.86%  " },
    @{ R = 17; D = "0.802";      E = "  -3.25%  " },
    @{ R = 18; D = "46.650.96";  E = "  -0.07%  " },
    @{ R = 21; D = "5.86";       E = "  -5.24%  " },
    @{ R = 22; D = "66.05";      E = "  -0.93%  " },
    @{ R = 23; D = "248.52";     E = "  +0.87%  " },
    @{ R = 24; D = "2.80";       E = "  -5.52%  " },
    @{ R = 25; D = "0.999";      E = "  -0.05%  " },
    @{ R = 26; D = $null;        E = "  -5.32%  " },
    @{ R = 27; D = "41.51";      E = "  -0.99%  " },
    @{ R = 28; D = $null;        E = "  -3.18%  " },
    @{ R = 29; D = "9.64";       E = "  -2.08%  " },
    @{ R = 30; D = "20.21";      E = "  +0.43%  " },
    @{ R = 31; D = $null;        E = "  +7.07%  " },
    @{ R = 32; D = "3.37";       E = "  +11.70%  " },
    @{ R = 33; D = "147.03";     E = "  -2.98%  " },
    @{ R = 34; D = "5.38";       E = "  -5.73%  " },
    @{ R = 35; D = "0.0772";     E = "  -4.44%  " },
    @{ R = 36; D = "0.114";      E = "  +7.65%  " },
    @{ R = 37; D = $null;        E = "  -2.52%  " },
    @{ R = 38; D = "15.86";      E = "  +14.55%  " },
    @{ R = 39; D = "1.69";       E = "  -5.93%  " },
    @{ R = 40; D = "3.90";       E = "  -2.94%  " },
    @{ R = 41; D = $null;        E = "  -6.25%  " },
    @{ R = 42; D = "3.12";       E = "  -8.64%  " },
    @{ R = 43; D = "0.998";      E = "  -0.17%  " },
    @{ R = 44; D = "93.82";      E = "  +15.96%  " },
    @{ R = 45; D = "1.793.88";   E = "  +0.06%  " },
    @{ R = 46; D = "1.87";       E = "  -4.28%  " },
    @{ R = 47; D = "71.40";      E = "  -2.94%  " },
    @{ R = 48; D = "0.186";      E = "  -6.38%  " },
    @{ R = 49; D = "4.82";       E = "  -1.63%  " },
    @{ R = 50; D = "7.95";       E = "  -0.08%  " },
    @{ R = 51; D = "95.29";      E = "  -2.75%  " }
)

foreach ($row in $rows) {
    if ($null -ne $row.D) {
        Set-TextValue "D$($row.R)" $row.D
    }
    Set-PlainValue "E$($row.R)" $row.E
}

# Rows 19 and 20 swapped position in the coin ranking: ShibaInu moved above
# InternetComputer(DFINITY). Update name, link, price and volume for both.
Set-PlainValue "B19" "ShibaInu"
Set-PlainValue "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue  "D19" "0.0₃0985"
Set-PlainValue "E19" "  +4.26%  "

Set-PlainValue "B20" "InternetComputer(DFINITY)"
Set-PlainValue "C20" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue  "D20" "12.55"
Set-PlainValue "E20" "  -4.53%  "
